$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 341 ("「ヨユウなきユウヨ」") was removed from the post list entirely;
# every subsequent row shifts up by one.
$ws.Rows.Item(341).Delete()
